$wb = $excel.ActiveWorkbook

# Citywide Totals
$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("D2").Value = 45
$ws.Range("I2").Value = 50
$ws.Range("E3").Value = 73
$ws.Range("C6").Value = 233
$ws.Range("E6").Value = 220
$ws.Range("F6").Value = 255
$ws.Range("G6").Value = 238
$ws.Range("H6").Value = 204
$ws.Range("I6").Value = 272
$ws.Range("C7").Value = 318
$ws.Range("D7").Value = 339
$ws.Range("E7").Value = 334
$ws.Range("F7").Value = 363
$ws.Range("G7").Value = 347
$ws.Range("H7").Value = 320
$ws.Range("I7").Value = 429

# Garfield Park
$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("E6").Value = 16
$ws.Range("H6").Value = 20
$ws.Range("E7").Value = 24
$ws.Range("H7").Value = 25

# Grand Crossing
$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("E4").Value = 13
$ws.Range("E5").Value = 16

# Uptown
$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("D2").Value = 1
$ws.Range("D6").Value = 2

# By Neighborhood
$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("E7").Value = 21
$ws.Range("E30").Value = 24
$ws.Range("H30").Value = 25
$ws.Range("G33").Value = 2
$ws.Range("E34").Value = 16
$ws.Range("I45").Value = 10
$ws.Range("C52").Value = 3
$ws.Range("I60").Value = 9
$ws.Range("D68").Value = 6
$ws.Range("F74").Value = 8
$ws.Range("D83").Value = 2
$ws.Range("C95").Value = 318
$ws.Range("D95").Value = 339
$ws.Range("E95").Value = 334
$ws.Range("F95").Value = 363
$ws.Range("G95").Value = 347
$ws.Range("H95").Value = 320
$ws.Range("I95").Value = 429

# Grand Boulevard
$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("G4").Value = 2
$ws.Range("G5").Value = 2

# Near South Side
$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("I2").Value = 1
$ws.Range("I5").Value = 9

# Roseland
$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("F4").Value = 2
$ws.Range("F5").Value = 8

# Lower West Side
$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("C4").Value = 3
$ws.Range("C5").Value = 3

# Lake View
$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("H5").Value = 9
$ws.Range("H6").Value = 10

# Austin
$ws = $wb.Worksheets.Item("Austin")
$ws.Range("E3").Value = 5
$ws.Range("E6").Value = 21

# Old Town
$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("C2").Value = 1
$ws.Range("C5").Value = 6
